$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new cat entry "Rita" as row 7
$ws.Range("A7").Value = "Rita"
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = "DSH"
$ws.Range("D7").Value = "Calico"
$ws.Range("E7").Value = "Social, vocal, sweet"

# Move selection to E9 (matches saved sheet view state)
$ws.Range("E9").Select() | Out-Null
